$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1 Kristyn Sergio")

# Update the Expertise (column G) values that changed between the
# before/after snapshots of this sheet.
$ws.Range("G3").Value  = "L"
$ws.Range("G5").Value  = "L"
$ws.Range("G6").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("G8").Value  = "L"
$ws.Range("G9").ClearContents()
$ws.Range("G11").Value = "M"
$ws.Range("G14").ClearContents()
$ws.Range("G15").Value = "L"
$ws.Range("G16").Value = "M"
$ws.Range("G17").Value = "L"
$ws.Range("G20").Value = "M"
$ws.Range("G25").Value = "L"
$ws.Range("G26").Value = "L"

# Update the active selection on the sheet (bottom-left pane selection
# moved from D23 to A10:K10).
$ws.Activate()
$ws.Range("A10:K10").Select()
